# TC for URL added
# Populates the "Actual Result" (J) and "Status" (K) columns for the
# URL_Test_Cases sheet, formats the new Status cells with a green fill,
# widens columns E and J, resizes a few rows to fit the new content, and
# updates the sheet's selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URL_Test_Cases")
$ws.Activate()

# ---------------------------------------------------------------------
# Column J (Actual Result) values - same wording style as the existing
# Expected Result column, plus a "Pass" Status in column K for each row.
# ---------------------------------------------------------------------

$ws.Range("J4").Value = "Application homepage loads successfully"
$ws.Range("K4").Value = "Pass"

$ws.Range("J5").Value = "Application launched successfully in Chrome, Firefox, and Edge"
$ws.Range("K5").Value = "Pass"

$ws.Range("J6").Value = "Browser displayed error`n page for invalid URL"
$ws.Range("J6").WrapText = $true
$ws.Range("K6").Value = "Pass"

$ws.Range("J7").Value = "Homepage loaded `nwithin 2.5 seconds"
$ws.Range("J7").WrapText = $true
$ws.Range("K7").Value = "Pass"

$ws.Range("J8").Value = "Application automatically `nredirected from HTTP to HTTPS"
$ws.Range("J8").WrapText = $true
$ws.Range("K8").Value = "Pass"

# ---------------------------------------------------------------------
# Style the new Status cells: green fill, vertical-centered wrapped text
# (order matters for how the style table gets built: wrap, then valign,
# then fill last, so no orphan intermediate styles are produced)
# ---------------------------------------------------------------------

$statusRange = $ws.Range("K4:K8")
$statusRange.WrapText = $true
$statusRange.VerticalAlignment = -4108  # xlCenter
$statusRange.Interior.Color = 5287936   # RGB(0, 176, 80) -> FF00B050

# ---------------------------------------------------------------------
# Row heights for rows whose wrapped content now needs more space
# ---------------------------------------------------------------------

$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 60

# ---------------------------------------------------------------------
# Column widths - widen E (Test Case Description) and J (Actual Result)
# to fit the new/longer content; drop the old auto bestFit sizing.
# ---------------------------------------------------------------------

$ws.Columns.Item(5).ColumnWidth = 52.76
$ws.Columns.Item(10).ColumnWidth = 21.17

# ---------------------------------------------------------------------
# Update the view: scroll so column E is at the left edge and select J10
# ---------------------------------------------------------------------

$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("J10").Select()
